$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: update title (D4) and link (E4)
$ws.Range("D4").Value = "처음 만나는 AI수학 with 파이썬 책 리뷰 - 수포자, 비전공자를 위한 수학"
$ws.Range("E4").Value = "https://teddylee777.github.io/thoughts/book-recomm-02"

# Row 23: update title (D23, multi-line) and link (E23)
$ws.Range("D23").Value = "안녕하세요! `nNLP를 공부하시는 분들, 혹은 처음 접하는 분들께 도움이 될까 하여 올려봅니다!  `n빅데이터 연합동아리 투빅스에서 7주동안"
$ws.Range("E23").Value = "https://theonly1.tistory.com/2686"

# Row 39: update title (D39) and link (E39)
$ws.Range("D39").Value = "KiCad 실습 - 12to5v 스텝다운 회로 아트웍"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/KiCad-%EC%8B%A4%EC%8A%B5-12to5v-%EC%8A%A4%ED%85%9D%EB%8B%A4%EC%9A%B4-%ED%9A%8C%EB%A1%9C-%EC%95%84%ED%8A%B8%EC%9B%8D-1"

# Row 46: update title (D46) and link (E46)
$ws.Range("D46").Value = "갈색세포종 땀분비"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/366"
